# Update "想去人数" (column F) values on the 展览 and 全部类型 sheets.
$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1.xml)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 11
$wsExhibit.Range("F3").Value = 5519
$wsExhibit.Range("F5").Value = 234
$wsExhibit.Range("F6").Value = 403
$wsExhibit.Range("F8").Value = 136
$wsExhibit.Range("F9").Value = 4338
$wsExhibit.Range("F10").Value = 780
$wsExhibit.Range("F11").Value = 804
$wsExhibit.Range("F12").Value = 40
$wsExhibit.Range("F18").Value = 118
$wsExhibit.Range("F19").Value = 603
$wsExhibit.Range("F20").Value = 19
$wsExhibit.Range("F22").Value = 1118
$wsExhibit.Range("F24").Value = 2746
$wsExhibit.Range("F25").Value = 438
$wsExhibit.Range("F26").Value = 284

# Sheet "全部类型" (sheet4.xml)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 11
$wsAll.Range("F3").Value = 5519
$wsAll.Range("F5").Value = 234
$wsAll.Range("F6").Value = 403
$wsAll.Range("F8").Value = 136
$wsAll.Range("F9").Value = 4338
$wsAll.Range("F10").Value = 780
$wsAll.Range("F11").Value = 804
$wsAll.Range("F12").Value = 40
$wsAll.Range("F18").Value = 118
$wsAll.Range("F19").Value = 603
$wsAll.Range("F20").Value = 19
$wsAll.Range("F23").Value = 1118
$wsAll.Range("F25").Value = 2746
$wsAll.Range("F26").Value = 438
$wsAll.Range("F27").Value = 284
